$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = 'Priority'
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 3
$ws.Range("G9").Value = 3
$ws.Range("G10").Value = 3
$ws.Range("B11").Value = 'MC'
$ws.Range("C11").Value = 2004
$ws.Range("D11").Value = 'Hull City'
$ws.Range("G11").Value = 1
$ws.Range("B12").Value = 'MCO'
$ws.Range("C12").Value = 2007
$ws.Range("D12").Value = 'Independiente del Valle'
$ws.Range("G12").Value = 1
$ws.Range("B13").Value = 'MC'
$ws.Range("C13").Value = 1999
$ws.Range("D13").Value = 'Botafogo-SP'
$ws.Range("E13").Value = 'Ghana'
$ws.Range("F13").Value = 'Segunda División Brasil'
$ws.Range("G13").Value = 2
$ws.Range("B14").Value = 'EI'
$ws.Range("C14").Value = 2005
$ws.Range("D14").Value = 'Santos'
$ws.Range("F14").Value = 'Primera División Brasil'
$ws.Range("G14").Value = 2
$ws.Range("B15").Value = 'DEF'
$ws.Range("C15").Value = 2006
$ws.Range("D15").Value = 'Independiente del Valle'
$ws.Range("E15").Value = 'Ecuador'
$ws.Range("F15").Value = 'Primera División Ecuador'
$ws.Range("G15").Value = 3
$ws.Range("B16").Value = 'DEL'
$ws.Range("C16").Value = 2006
$ws.Range("D16").Value = 'Barcelona SC'
$ws.Range("G16").Value = 2
$ws.Range("B17").Value = 'MC'
$ws.Range("C17").Value = 2004
$ws.Range("D17").Value = 'Pacos Ferreira'
$ws.Range("E17").Value = 'Ecuador'
$ws.Range("F17").Value = 'Segunda División Portugal'
$ws.Range("G17").Value = 2
$ws.Range("B18").Value = 'MCO'
$ws.Range("C18").Value = 2006
$ws.Range("D18").Value = 'Independiente del Valle'
$ws.Range("E18").Value = 'Ecuador'
$ws.Range("F18").Value = 'Primera División Ecuador'
$ws.Range("G18").Value = 3
$ws.Range("B19").Value = 'EI'
$ws.Range("C19").Value = 2006
$ws.Range("D19").Value = 'Independiente del Valle'
$ws.Range("E19").Value = 'Ecuador'
$ws.Range("F19").Value = 'Primera División Ecuador'
$ws.Range("G19").Value = 3
$ws.Range("B20").Value = 'MC'
$ws.Range("C20").Value = 2001
$ws.Range("D20").Value = 'Chelsea'
$ws.Range("G20").Value = 1
$ws.Range("B21").Value = 'MC'
$ws.Range("C21").Value = 2008
$ws.Range("D21").Value = 'Independiente del Valle'
$ws.Range("E21").Value = 'Ecuador'
$ws.Range("F21").Value = 'Primera División Ecuador'
$ws.Range("G21").Value = 3
$ws.Range("B22").Value = 'DEF'
$ws.Range("C22").Value = 1997
$ws.Range("D22").Value = 'Ethnikos'
$ws.Range("E22").Value = 'Ghana'
$ws.Range("F22").Value = 'Primera División Chipre'
$ws.Range("G22").Value = 2
$ws.Range("B23").Value = 'ED'
$ws.Range("C23").Value = 2000
$ws.Range("D23").Value = 'Flamengo'
$ws.Range("G23").Value = 1

$ws.Range("G24").Select()
try {
    $excel.ActiveWindow.ScrollRow = 6
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
